$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Section 1: "Preliminary BOM"
# ---------------------------------------------------------------------------
$ws.Range("A7").Value = "Preliminary BOM"
$ws.Range("H7").Value = "Prelim. BOM Total"

$ws.Range("A8").Value = "Manufacturer"
$ws.Range("B8").Value = "Manufacturer SKU"
$ws.Range("C8").Value = "Digikey SKU"
$ws.Range("D8").Value = "Quantity"
$ws.Range("E8").Value = "Price"
$ws.Range("F8").Value = "Subtotal"
$ws.Range("H8").Formula = "=SUM(F:F)"

$ws.Range("A9").Value = "Mill-Max"
$ws.Range("B9").Value = "2906-4-15-20-75-14-11-0"
$ws.Range("C9").Value = "54-2906-4-15-20-75-14-11-0-ND"
$ws.Range("D9").Value = 14
$ws.Range("E9").Value = 0.64
$ws.Range("F9").Formula = "=E9*D9"

$ws.Range("A10").Value = "GCT"
$ws.Range("B10").Value = "USB4056-03-A"
$ws.Range("C10").Value = "2073-USB4056-03-ATR-ND"
$ws.Range("D10").Value = 2
$ws.Range("E10").Value = 0.97
$ws.Range("F10").Formula = "=E10*D10"

$ws.Range("A11").Value = "GuHua"
$ws.Range("B11").Value = "Aliexpress Heat Set Inserts"
$ws.Range("C11").Value = "n/a"
$ws.Range("D11").Value = 8
$ws.Range("E11").Value = 0.047
$ws.Range("F11").Formula = "=E11*D11"

$ws.Range("A12").Value = "K&J Magnetics"
$ws.Range("B12").Value = "D73-N52"
$ws.Range("C12").Value = "n/a"
$ws.Range("D12").Value = 4
$ws.Range("E12").Formula = "=1.35"
$ws.Range("F12").Formula = "=E12*D12"

$ws.Range("A13").Value = "PCBWay"
$ws.Range("B13").Value = "n/a"
$ws.Range("C13").Value = "n/a"
$ws.Range("D13").Value = 2
$ws.Range("E13").Value = 10
$ws.Range("F13").Formula = "=E13*D13"

# ---------------------------------------------------------------------------
# Section 2: "Printed Parts"
# ---------------------------------------------------------------------------
$ws.Range("A14").Value = "Printed Parts"

$ws.Range("A15").Value = "Filament"
$ws.Range("B15").Value = "Part"
$ws.Range("C15").Value = "Estimated Amount (g)"
$ws.Range("D15").Value = "Quantity"
$ws.Range("E15").Value = "Price"
$ws.Range("F15").Value = "Subtotal"

$ws.Range("A16").Value = "Gen. ASA"
$ws.Range("B16").Value = "Female Half"
$ws.Range("C16").Value = 10
$ws.Range("D16").Value = 1
$ws.Range("E16").Formula = "=(25/2000)*C16"
$ws.Range("F16").Formula = "=D16*E16"

$ws.Range("A17").Value = "Gen. ASA"
$ws.Range("B17").Value = "Male Half"
$ws.Range("C17").Value = 10
$ws.Range("D17").Value = 1
$ws.Range("E17").Formula = "=(25/2000)*C17"
$ws.Range("F17").Formula = "=D17*E17"

# ---------------------------------------------------------------------------
# Section 3: "Equipment Tracking"
# ---------------------------------------------------------------------------
$ws.Range("A20").Value = "Equipment Tracking"

$ws.Range("A21").Value = "Name"
$ws.Range("B21").Value = "Price"
$ws.Range("C21").Value = "Total"

$ws.Range("A22").Value = "Soldering Plate"
$ws.Range("B22").Value = 24.93
$ws.Range("C22").Formula = "=SUM(B22:B32)"

$ws.Range("A23").Value = "Solder Paste"
$ws.Range("B23").Value = 7.1

$ws.Range("A24").Value = "Flux Paste"
$ws.Range("B24").Value = 4.1

$ws.Range("B25").Value = $null

Write-Host "values set"
